$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("D2").Value = 15235
$ws.Range("E2").Value = 398
$ws.Range("F2").Value = 398
$ws.Range("G2").Value = 1748
$ws.Range("H2").Value = 1058
$ws.Range("I2").Value = 1058
$ws.Range("K2").Value = 15044
$ws.Range("L2").Value = 12035
$ws.Range("M2").Value = 3009
$ws.Range("N2").Value = 3009
$ws.Range("P2").Value = 1716
$ws.Range("Q2").Value = -1894
$ws.Range("R2").Value = 1207
$ws.Range("S2").Value = -416
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = -1900
$ws.Range("V2").Value = 4356
$ws.Range("W2").Value = 2.61
$ws.Range("X2").Value = 6.94
$ws.Range("Y2").Value = 46.6
$ws.Range("Z2").Value = 6.47
$ws.Range("AA2").Value = 399.95
$ws.Range("AB2").Value = 63.2
$ws.Range("AC2").Value = 3194
$ws.Range("AD2").Value = 7.11
$ws.Range("AE2").Value = 8776
$ws.Range("AF2").Value = 2.59
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 34017685

# Row 3
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("D3").Value = 15310
$ws.Range("E3").Value = 208
$ws.Range("F3").Value = 208
$ws.Range("G3").Value = -147
$ws.Range("H3").Value = -61
$ws.Range("I3").Value = -61
$ws.Range("K3").Value = 13100
$ws.Range("L3").Value = 10199
$ws.Range("M3").Value = 2902
$ws.Range("N3").Value = 2902
$ws.Range("P3").Value = 1755
$ws.Range("Q3").Value = 182
$ws.Range("R3").Value = 446
$ws.Range("S3").Value = -1186
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = 175
$ws.Range("V3").Value = 3345
$ws.Range("W3").Value = 1.36
$ws.Range("X3").Value = -0.4
$ws.Range("Y3").Value = -2.07
$ws.Range("Z3").Value = -0.43
$ws.Range("AA3").Value = 351.47
$ws.Range("AB3").Value = 59.07
$ws.Range("AC3").Value = -176
$ws.Range("AD3").Value = -84.29
$ws.Range("AE3").Value = 8272
$ws.Range("AF3").Value = 1.79
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 34806923

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("D4").Value = 13537
$ws.Range("E4").Value = 418
$ws.Range("F4").Value = 418
$ws.Range("G4").Value = 781
$ws.Range("H4").Value = 360
$ws.Range("I4").Value = 360
$ws.Range("K4").Value = 12746
$ws.Range("L4").Value = 9633
$ws.Range("M4").Value = 3113
$ws.Range("N4").Value = 3113
$ws.Range("P4").Value = 1769
$ws.Range("Q4").Value = 191
$ws.Range("R4").Value = 222
$ws.Range("S4").Value = -875
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = 186
$ws.Range("V4").Value = 2689
$ws.Range("W4").Value = 3.09
$ws.Range("X4").Value = 2.66
$ws.Range("Y4").Value = 11.98
$ws.Range("Z4").Value = 2.79
$ws.Range("AA4").Value = 309.5
$ws.Range("AB4").Value = 77.37
$ws.Range("AC4").Value = 1021
$ws.Range("AD4").Value = 9.48
$ws.Range("AE4").Value = 8802
$ws.Range("AF4").Value = 1.1
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 3.1
$ws.Range("AI4").Value = 29.48
$ws.Range("AJ4").Value = 35090993

# Row 5
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("D5").Value = 12979
$ws.Range("E5").Value = 311
$ws.Range("F5").Value = 311
$ws.Range("G5").Value = 998
$ws.Range("H5").Value = 953
$ws.Range("I5").Value = 953
$ws.Range("K5").Value = 12849
$ws.Range("L5").Value = 8922
$ws.Range("M5").Value = 3927
$ws.Range("N5").Value = 3927
$ws.Range("P5").Value = 1807
$ws.Range("Q5").Value = 759
$ws.Range("R5").Value = -103
$ws.Range("S5").Value = -826
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 754
$ws.Range("V5").Value = 2172
$ws.Range("W5").Value = 2.39
$ws.Range("X5").Value = 7.34
$ws.Range("Y5").Value = 27.07
$ws.Range("Z5").Value = 7.45
$ws.Range("AA5").Value = 227.2
$ws.Range("AB5").Value = 124.08
$ws.Range("AC5").Value = 2652
$ws.Range("AD5").Value = 3.59
$ws.Range("AE5").Value = 11115
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 5.26
$ws.Range("AI5").Value = 18.56
$ws.Range("AJ5").Value = 35849527

# Row 6
$ws.Range("D6").Value = 13767
$ws.Range("E6").Value = 423
$ws.Range("F6").Value = 423
$ws.Range("G6").Value = 153
$ws.Range("H6").Value = -5
$ws.Range("I6").Value = -5
$ws.Range("K6").Value = 11822
$ws.Range("L6").Value = 8259
$ws.Range("M6").Value = 3563
$ws.Range("N6").Value = 3563
$ws.Range("P6").Value = 1816
$ws.Range("Q6").Value = 557
$ws.Range("R6").Value = -178
$ws.Range("S6").Value = -507
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 553
$ws.Range("V6").Value = 1923
$ws.Range("W6").Value = 3.07
$ws.Range("X6").Value = -0.03
$ws.Range("Y6").Value = -0.13
$ws.Range("Z6").Value = -0.04
$ws.Range("AA6").Value = 231.8
$ws.Range("AB6").Value = 128.67
$ws.Range("AC6").Value = -13
$ws.Range("AD6").Value = -907.86
$ws.Range("AE6").Value = 10103
$ws.Range("AF6").Value = 1.17
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 4.24
$ws.Range("AI6").Value = -3750.33
$ws.Range("AJ6").Value = 36035349

# Row 7
$ws.Range("D7").Value = 16272
$ws.Range("E7").Value = 558
$ws.Range("G7").Value = 299
$ws.Range("H7").Value = 118
$ws.Range("I7").Value = 248
$ws.Range("K7").Value = 13113
$ws.Range("L7").Value = 9265
$ws.Range("M7").Value = 3850
$ws.Range("N7").Value = 3880
$ws.Range("P7").Value = 1819
$ws.Range("Q7").Value = 1013
$ws.Range("R7").Value = -81
$ws.Range("S7").Value = -380
$ws.Range("T7").Value = 304
$ws.Range("U7").Value = -76
$ws.Range("W7").Value = 3.43
$ws.Range("X7").Value = 0.73
$ws.Range("Y7").Value = 6.66
$ws.Range("Z7").Value = 0.95
$ws.Range("AA7").Value = 240.68
$ws.Range("AC7").Value = 682
$ws.Range("AD7").Value = 13.88
$ws.Range("AE7").Value = 11000
$ws.Range("AF7").Value = 0.86
$ws.Range("AG7").Value = 540
$ws.Range("AH7").Value = 5.7
$ws.Range("AI7").Value = 78.52

# Row 8
$ws.Range("D8").Value = 17490
$ws.Range("E8").Value = 760
$ws.Range("G8").Value = 798
$ws.Range("H8").Value = 675
$ws.Range("I8").Value = 644
$ws.Range("K8").Value = 13685
$ws.Range("L8").Value = 9364
$ws.Range("M8").Value = 4324
$ws.Range("N8").Value = 4412
$ws.Range("P8").Value = 1819
$ws.Range("Q8").Value = 886
$ws.Range("R8").Value = -154
$ws.Range("S8").Value = -374
$ws.Range("T8").Value = 15
$ws.Range("U8").Value = 337
$ws.Range("W8").Value = 4.35
$ws.Range("X8").Value = 3.86
$ws.Range("Y8").Value = 15.53
$ws.Range("Z8").Value = 5.04
$ws.Range("AA8").Value = 216.57
$ws.Range("AC8").Value = 1773
$ws.Range("AD8").Value = 5.34
$ws.Range("AE8").Value = 12510
$ws.Range("AF8").Value = 0.76
$ws.Range("AG8").Value = 570
$ws.Range("AH8").Value = 6.02
$ws.Range("AI8").Value = 31.89

# Row 9
$ws.Range("D9").Value = 18537
$ws.Range("E9").Value = 876
$ws.Range("G9").Value = 763
$ws.Range("H9").Value = 640
$ws.Range("I9").Value = 620
$ws.Range("K9").Value = 14400
$ws.Range("L9").Value = 9551
$ws.Range("M9").Value = 4846
$ws.Range("N9").Value = 4846
$ws.Range("P9").Value = 1819
$ws.Range("Q9").Value = 1084
$ws.Range("R9").Value = -230
$ws.Range("S9").Value = -437
$ws.Range("T9").Value = 27
$ws.Range("U9").Value = 719
$ws.Range("W9").Value = 4.72
$ws.Range("X9").Value = 3.45
$ws.Range("Y9").Value = 13.39
$ws.Range("Z9").Value = 4.56
$ws.Range("AA9").Value = 197.1
$ws.Range("AC9").Value = 1706
$ws.Range("AD9").Value = 5.55
$ws.Range("AE9").Value = 13742
$ws.Range("AF9").Value = 0.69
$ws.Range("AG9").Value = 550
$ws.Range("AH9").Value = 5.81
$ws.Range("AI9").Value = 31.98
